$wb = $excel.ActiveWorkbook

# ---- PIR ----
$ws = $wb.Worksheets.Item('PIR')
$ws.Range("A567:A579").NumberFormat = "@"
$ws.Cells.Item(567, 1).Value = '2026-02-06'
$ws.Cells.Item(567, 2).Value = '10:27:02'
$ws.Cells.Item(567, 3).Value = '10:00'
$ws.Cells.Item(567, 4).Value = 'Bathroom'
$ws.Cells.Item(567, 5).Value = 'No Motion'
$ws.Cells.Item(567, 6).Value = 'Inactive'
$ws.Cells.Item(568, 1).Value = '2026-02-06'
$ws.Cells.Item(568, 2).Value = '10:27:06'
$ws.Cells.Item(568, 3).Value = '10:00'
$ws.Cells.Item(568, 4).Value = 'Bathroom'
$ws.Cells.Item(568, 5).Value = 'No Motion'
$ws.Cells.Item(568, 6).Value = 'Inactive'
$ws.Cells.Item(569, 1).Value = '2026-02-06'
$ws.Cells.Item(569, 2).Value = '10:27:07'
$ws.Cells.Item(569, 3).Value = '10:00'
$ws.Cells.Item(569, 4).Value = 'Bathroom'
$ws.Cells.Item(569, 5).Value = 'No Motion'
$ws.Cells.Item(569, 6).Value = 'Inactive'
$ws.Cells.Item(570, 1).Value = '2026-02-06'
$ws.Cells.Item(570, 2).Value = '10:27:12'
$ws.Cells.Item(570, 3).Value = '10:00'
$ws.Cells.Item(570, 4).Value = 'Bathroom'
$ws.Cells.Item(570, 5).Value = 'No Motion'
$ws.Cells.Item(570, 6).Value = 'Inactive'
$ws.Cells.Item(571, 1).Value = '2026-02-06'
$ws.Cells.Item(571, 2).Value = '10:27:17'
$ws.Cells.Item(571, 3).Value = '10:00'
$ws.Cells.Item(571, 4).Value = 'Bathroom'
$ws.Cells.Item(571, 5).Value = 'No Motion'
$ws.Cells.Item(571, 6).Value = 'Inactive'
$ws.Cells.Item(572, 1).Value = '2026-02-06'
$ws.Cells.Item(572, 2).Value = '10:27:22'
$ws.Cells.Item(572, 3).Value = '10:00'
$ws.Cells.Item(572, 4).Value = 'Bathroom'
$ws.Cells.Item(572, 5).Value = 'No Motion'
$ws.Cells.Item(572, 6).Value = 'Inactive'
$ws.Cells.Item(573, 1).Value = '2026-02-06'
$ws.Cells.Item(573, 2).Value = '10:27:27'
$ws.Cells.Item(573, 3).Value = '10:00'
$ws.Cells.Item(573, 4).Value = 'Bathroom'
$ws.Cells.Item(573, 5).Value = 'No Motion'
$ws.Cells.Item(573, 6).Value = 'Inactive'
$ws.Cells.Item(574, 1).Value = '2026-02-06'
$ws.Cells.Item(574, 2).Value = '10:27:32'
$ws.Cells.Item(574, 3).Value = '10:00'
$ws.Cells.Item(574, 4).Value = 'Bathroom'
$ws.Cells.Item(574, 5).Value = 'No Motion'
$ws.Cells.Item(574, 6).Value = 'Inactive'
$ws.Cells.Item(575, 1).Value = '2026-02-06'
$ws.Cells.Item(575, 2).Value = '10:27:37'
$ws.Cells.Item(575, 3).Value = '10:00'
$ws.Cells.Item(575, 4).Value = 'Bathroom'
$ws.Cells.Item(575, 5).Value = 'No Motion'
$ws.Cells.Item(575, 6).Value = 'Inactive'
$ws.Cells.Item(576, 1).Value = '2026-02-06'
$ws.Cells.Item(576, 2).Value = '10:27:42'
$ws.Cells.Item(576, 3).Value = '10:00'
$ws.Cells.Item(576, 4).Value = 'Bathroom'
$ws.Cells.Item(576, 5).Value = 'No Motion'
$ws.Cells.Item(576, 6).Value = 'Inactive'
$ws.Cells.Item(577, 1).Value = '2026-02-06'
$ws.Cells.Item(577, 2).Value = '10:27:48'
$ws.Cells.Item(577, 3).Value = '10:00'
$ws.Cells.Item(577, 4).Value = 'Bathroom'
$ws.Cells.Item(577, 5).Value = 'No Motion'
$ws.Cells.Item(577, 6).Value = 'Inactive'
$ws.Cells.Item(578, 1).Value = '2026-02-06'
$ws.Cells.Item(578, 2).Value = '10:27:53'
$ws.Cells.Item(578, 3).Value = '10:00'
$ws.Cells.Item(578, 4).Value = 'Bathroom'
$ws.Cells.Item(578, 5).Value = 'No Motion'
$ws.Cells.Item(578, 6).Value = 'Inactive'
$ws.Cells.Item(579, 1).Value = '2026-02-06'
$ws.Cells.Item(579, 2).Value = '10:27:58'
$ws.Cells.Item(579, 3).Value = '10:00'
$ws.Cells.Item(579, 4).Value = 'Bathroom'
$ws.Cells.Item(579, 5).Value = 'No Motion'
$ws.Cells.Item(579, 6).Value = 'Inactive'
$ws.Range("A567:A579").ClearFormats()

# ---- Humidity ----
$ws = $wb.Worksheets.Item('Humidity')
$ws.Range("A402:A411").NumberFormat = "@"
$ws.Range("E402:E411").NumberFormat = "@"
$ws.Cells.Item(402, 1).Value = '2026-02-06'
$ws.Cells.Item(402, 2).Value = '10:27:03'
$ws.Cells.Item(402, 3).Value = '10:00'
$ws.Cells.Item(402, 4).Value = 'Bathroom'
$ws.Cells.Item(402, 5).Value = '67.6%'
$ws.Cells.Item(402, 6).Value = 'Active'
$ws.Cells.Item(403, 1).Value = '2026-02-06'
$ws.Cells.Item(403, 2).Value = '10:27:09'
$ws.Cells.Item(403, 3).Value = '10:00'
$ws.Cells.Item(403, 4).Value = 'Bathroom'
$ws.Cells.Item(403, 5).Value = '67.7%'
$ws.Cells.Item(403, 6).Value = 'Active'
$ws.Cells.Item(404, 1).Value = '2026-02-06'
$ws.Cells.Item(404, 2).Value = '10:27:19'
$ws.Cells.Item(404, 3).Value = '10:00'
$ws.Cells.Item(404, 4).Value = 'Bathroom'
$ws.Cells.Item(404, 5).Value = '67.8%'
$ws.Cells.Item(404, 6).Value = 'Active'
$ws.Cells.Item(405, 1).Value = '2026-02-06'
$ws.Cells.Item(405, 2).Value = '10:27:24'
$ws.Cells.Item(405, 3).Value = '10:00'
$ws.Cells.Item(405, 4).Value = 'Bathroom'
$ws.Cells.Item(405, 5).Value = '67.7%'
$ws.Cells.Item(405, 6).Value = 'Active'
$ws.Cells.Item(406, 1).Value = '2026-02-06'
$ws.Cells.Item(406, 2).Value = '10:27:29'
$ws.Cells.Item(406, 3).Value = '10:00'
$ws.Cells.Item(406, 4).Value = 'Bathroom'
$ws.Cells.Item(406, 5).Value = '67.8%'
$ws.Cells.Item(406, 6).Value = 'Active'
$ws.Cells.Item(407, 1).Value = '2026-02-06'
$ws.Cells.Item(407, 2).Value = '10:27:34'
$ws.Cells.Item(407, 3).Value = '10:00'
$ws.Cells.Item(407, 4).Value = 'Bathroom'
$ws.Cells.Item(407, 5).Value = '67.8%'
$ws.Cells.Item(407, 6).Value = 'Active'
$ws.Cells.Item(408, 1).Value = '2026-02-06'
$ws.Cells.Item(408, 2).Value = '10:27:39'
$ws.Cells.Item(408, 3).Value = '10:00'
$ws.Cells.Item(408, 4).Value = 'Bathroom'
$ws.Cells.Item(408, 5).Value = '67.7%'
$ws.Cells.Item(408, 6).Value = 'Active'
$ws.Cells.Item(409, 1).Value = '2026-02-06'
$ws.Cells.Item(409, 2).Value = '10:27:44'
$ws.Cells.Item(409, 3).Value = '10:00'
$ws.Cells.Item(409, 4).Value = 'Bathroom'
$ws.Cells.Item(409, 5).Value = '67.7%'
$ws.Cells.Item(409, 6).Value = 'Active'
$ws.Cells.Item(410, 1).Value = '2026-02-06'
$ws.Cells.Item(410, 2).Value = '10:27:54'
$ws.Cells.Item(410, 3).Value = '10:00'
$ws.Cells.Item(410, 4).Value = 'Bathroom'
$ws.Cells.Item(410, 5).Value = '67.5%'
$ws.Cells.Item(410, 6).Value = 'Active'
$ws.Cells.Item(411, 1).Value = '2026-02-06'
$ws.Cells.Item(411, 2).Value = '10:27:59'
$ws.Cells.Item(411, 3).Value = '10:00'
$ws.Cells.Item(411, 4).Value = 'Bathroom'
$ws.Cells.Item(411, 5).Value = '67.5%'
$ws.Cells.Item(411, 6).Value = 'Active'
$ws.Range("A402:A411").ClearFormats()
$ws.Range("E402:E411").ClearFormats()

# ---- Temperature ----
$ws = $wb.Worksheets.Item('Temperature')
$ws.Range("A402:A411").NumberFormat = "@"
$ws.Cells.Item(402, 1).Value = '2026-02-06'
$ws.Cells.Item(402, 2).Value = '10:27:04'
$ws.Cells.Item(402, 3).Value = '10:00'
$ws.Cells.Item(402, 4).Value = 'Bathroom'
$ws.Cells.Item(402, 5).Value = '28.4C'
$ws.Cells.Item(402, 6).Value = 'Active'
$ws.Cells.Item(403, 1).Value = '2026-02-06'
$ws.Cells.Item(403, 2).Value = '10:27:10'
$ws.Cells.Item(403, 3).Value = '10:00'
$ws.Cells.Item(403, 4).Value = 'Bathroom'
$ws.Cells.Item(403, 5).Value = '28.5C'
$ws.Cells.Item(403, 6).Value = 'Active'
$ws.Cells.Item(404, 1).Value = '2026-02-06'
$ws.Cells.Item(404, 2).Value = '10:27:20'
$ws.Cells.Item(404, 3).Value = '10:00'
$ws.Cells.Item(404, 4).Value = 'Bathroom'
$ws.Cells.Item(404, 5).Value = '28.4C'
$ws.Cells.Item(404, 6).Value = 'Active'
$ws.Cells.Item(405, 1).Value = '2026-02-06'
$ws.Cells.Item(405, 2).Value = '10:27:25'
$ws.Cells.Item(405, 3).Value = '10:00'
$ws.Cells.Item(405, 4).Value = 'Bathroom'
$ws.Cells.Item(405, 5).Value = '28.3C'
$ws.Cells.Item(405, 6).Value = 'Active'
$ws.Cells.Item(406, 1).Value = '2026-02-06'
$ws.Cells.Item(406, 2).Value = '10:27:30'
$ws.Cells.Item(406, 3).Value = '10:00'
$ws.Cells.Item(406, 4).Value = 'Bathroom'
$ws.Cells.Item(406, 5).Value = '28.3C'
$ws.Cells.Item(406, 6).Value = 'Active'
$ws.Cells.Item(407, 1).Value = '2026-02-06'
$ws.Cells.Item(407, 2).Value = '10:27:35'
$ws.Cells.Item(407, 3).Value = '10:00'
$ws.Cells.Item(407, 4).Value = 'Bathroom'
$ws.Cells.Item(407, 5).Value = '28.4C'
$ws.Cells.Item(407, 6).Value = 'Active'
$ws.Cells.Item(408, 1).Value = '2026-02-06'
$ws.Cells.Item(408, 2).Value = '10:27:40'
$ws.Cells.Item(408, 3).Value = '10:00'
$ws.Cells.Item(408, 4).Value = 'Bathroom'
$ws.Cells.Item(408, 5).Value = '28.3C'
$ws.Cells.Item(408, 6).Value = 'Active'
$ws.Cells.Item(409, 1).Value = '2026-02-06'
$ws.Cells.Item(409, 2).Value = '10:27:45'
$ws.Cells.Item(409, 3).Value = '10:00'
$ws.Cells.Item(409, 4).Value = 'Bathroom'
$ws.Cells.Item(409, 5).Value = '28.4C'
$ws.Cells.Item(409, 6).Value = 'Active'
$ws.Cells.Item(410, 1).Value = '2026-02-06'
$ws.Cells.Item(410, 2).Value = '10:27:55'
$ws.Cells.Item(410, 3).Value = '10:00'
$ws.Cells.Item(410, 4).Value = 'Bathroom'
$ws.Cells.Item(410, 5).Value = '28.3C'
$ws.Cells.Item(410, 6).Value = 'Active'
$ws.Cells.Item(411, 1).Value = '2026-02-06'
$ws.Cells.Item(411, 2).Value = '10:28:00'
$ws.Cells.Item(411, 3).Value = '10:00'
$ws.Cells.Item(411, 4).Value = 'Bathroom'
$ws.Cells.Item(411, 5).Value = '28.3C'
$ws.Cells.Item(411, 6).Value = 'Active'
$ws.Range("A402:A411").ClearFormats()

